$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the formula in B1 (=160*10) with a plain literal value
$ws.Range("B1").Value = 400

# Move the active selection from B3 to B1
$ws.Range("B1").Select()
